$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Rename Sheet2 -> "Database testing" ---
$ws2.Name = "Database testing"

# --- Populate header row (row 1) on the "Database testing" sheet, same as Sheet1 ---
$ws2.Range("A1").Value = "name"
$ws2.Range("B1").Value = "surname"
$ws2.Range("C1").Value = "language"
$ws2.Range("D1").Value = "add1"
$ws2.Range("E1").Value = "add2"
$ws2.Range("F1").Value = "city"
$ws2.Range("G1").Value = "state"
$ws2.Range("H1").Value = "zip"
$ws2.Range("I1").Value = "country"
$ws2.Range("J1").Value = "gender"
$ws2.Range("K1").NumberFormat = "@"
$ws2.Range("K1").Value = "birthdate"
$ws2.Range("L1").Value = "phone"
$ws2.Range("M1").Value = "faxNumber"
$ws2.Range("N1").Value = "Mob"
$ws2.Range("O1").Value = "email"
$ws2.Range("P1").Value = "web"
$ws2.Range("Q1").Value = "vat"
$ws2.Range("R1").Value = "tax"

# --- Populate data row (row 2) - mirrors Sheet1 row 2, with a new "ajay3" record ---
$ws2.Range("A2").Value = "ajay3"
$ws2.Range("B2").Value = "ghodake"
$ws2.Range("C2").Value = "English"
$ws2.Range("D2").Value = "NewJersey"
$ws2.Range("E2").Value = "xyz"
$ws2.Range("F2").Value = "Edison "
$ws2.Range("G2").Value = "NJ"
$ws2.Range("H2").NumberFormat = "00000"
$ws2.Range("H2").Value = 8827
$ws2.Range("I2").Value = "United States"
$ws2.Range("J2").Value = "Male"
$ws2.Range("K2").NumberFormat = "@"
$ws2.Range("K2").Value = "23/08/2020"
$ws2.Range("L2").Value = "567-999-7456"
$ws2.Range("M2").Value = 45678
$ws2.Range("N2").Value = 76543
$ws2.Range("O2").Value = "a@b.com"
$ws2.Range("P2").Value = "www.xyz"
$ws2.Range("Q2").Value = 434
$ws2.Range("R2").Value = 3434

# --- Hyperlinks for P2 (web) and O2 (email) mirroring Sheet1 (P2 added first) ---
$ws2.Hyperlinks.Add($ws2.Range("P2"), "http://www.xyz/")
$ws2.Hyperlinks.Add($ws2.Range("O2"), "mailto:a@b.com")
# Restore the shared "Hyperlink" cell style (Add() bumps a fresh style variant otherwise)
$ws2.Range("O2").Style = "Hyperlink"
$ws2.Range("P2").Style = "Hyperlink"

# Row 2 is slightly taller on Sheet1 (15.75) - match it here too
$ws2.Rows(2).RowHeight = 15.75

# --- View state: Sheet1 loses the active selection/tab, "Database testing" becomes active ---
$ws1.Activate()
$ws1.Range("A1:XFD2").Select()

$ws2.Activate()
$ws2.Range("C5").Select()
$excel.ActiveWindow.Zoom = 145
